# Insert a new "Drone Applications" slide at position 3 (Title and Content layout).
$p = $ppt.ActivePresentation

$s = $p.Slides.Add(3, 2)

# --- Title placeholder ---
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Drone Applications"

# --- Content placeholder ---
$body = $s.Shapes.Item(2)
$tr = $body.TextFrame.TextRange
$tr.Text = "Crop Monitoring " + [char]0x2013 + " Agriculture`rTraffic Breach " + [char]0x2013 + " Surveillance`rMedical Transport System " + [char]0x2013 + " Healthcare`rDisaster Management " + [char]0x2013 + " Surveillance "

# Split the final paragraph into two runs: "...Surveillance " + "and Security"
[void]$tr.InsertAfter("and Security")

# Resize/position the content placeholder to match the authored layout tweak.
$body.TextFrame.AutoSize = 2
$body.Height = 3008612 / 12700
